$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.276.95"
$ws.Range("E2").Value = "  +2.16%  "

$ws.Range("D3").Value = "2.059.34"
$ws.Range("E3").Value = "  +3.53%  "

$ws.Range("E4").Value = "  -0.12%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.25"
$ws.Range("E5").Value = "  -0.80%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.611"
$ws.Range("E6").Value = "  +2.76%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.14"
$ws.Range("E7").Value = "  +6.06%  "

$ws.Range("E8").Value = "  -0.05%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.382"
$ws.Range("E9").Value = "  +2.95%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "58.69"
$ws.Range("E10").Value = "  +2.21%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0759"
$ws.Range("E11").Value = "  +1.70%  "

$ws.Range("E12").Value = "  +2.67%  "

$ws.Range("D13").Value = "2.366.24"
$ws.Range("E13").Value = "  +3.67%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.61"
$ws.Range("E14").Value = "  +3.49%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.05"
$ws.Range("E15").Value = "  +3.06%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.774"
$ws.Range("E16").Value = "  +2.00%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.16"
$ws.Range("E17").Value = "  +1.86%  "

$ws.Range("D18").Value = "2.093.36"
$ws.Range("E18").Value = "  +4.96%  "

$ws.Range("D19").Value = "37.503.49"
$ws.Range("E19").Value = "  +2.86%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.12"
$ws.Range("E20").Value = "  +16.54%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "68.84"
$ws.Range("E21").Value = "  +1.62%  "

$ws.Range("D22").Value = "0.0₃0812"
$ws.Range("E22").Value = "  +0.93%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "226.21"
$ws.Range("E23").Value = "  +2.16%  "

$ws.Range("E24").Value = "  +0.02%  "

$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.40"
$ws.Range("E25").Value = "  +1.44%  "

$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.42"
$ws.Range("E26").Value = "  +1.09%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "165.00"
$ws.Range("E27").Value = "  +1.78%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.50"
$ws.Range("E28").Value = "  +13.30%  "

$ws.Range("E29").Value = "  +2.01%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.14"
$ws.Range("E30").Value = "  +1.44%  "

$ws.Range("E31").Value = "  -1.81%  "

$ws.Range("E32").Value = "  +1.44%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.51"
$ws.Range("E33").Value = "  +3.35%  "

$ws.Range("E34").Value = "  +2.79%  "

$ws.Range("E35").Value = "  +8.97%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.55"
$ws.Range("E36").Value = "  +6.72%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.37"
$ws.Range("E37").Value = "  -1.82%  "

$ws.Range("E38").Value = "  +0.01%  "

$ws.Range("E39").Value = "  +0.93%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.84"
$ws.Range("E40").Value = "  +4.10%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0985"
$ws.Range("E41").Value = "  +4.75%  "

$ws.Range("E42").Value = "  -1.74%  "

$ws.Range("B43").Value = "FTXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.34"
$ws.Range("E43").Value = "  +19.93%  "

$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "1.456.73"
$ws.Range("E44").Value = "  +0.19%  "

$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "95.68"
$ws.Range("E45").Value = "  +7.44%  "

$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0210"
$ws.Range("E46").Value = "  +3.44%  "

$ws.Range("B47").Value = "TrustWalletToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.15"
$ws.Range("E47").Value = "  +4.66%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "15.84"
$ws.Range("E48").Value = "  +4.65%  "

$ws.Range("E49").Value = "  +3.40%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.26"
$ws.Range("E50").Value = "  +5.89%  "

$ws.Range("E51").Value = "  +2.07%  "
